$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.21"
$ws.Range("E2").Value = "'2.04%"
$ws.Range("G2").Value = "'22"
$ws.Range("D3").Value = "'29.59"
$ws.Range("E3").Value = "'3.36%"
$ws.Range("G3").Value = "'22"
$ws.Range("D4").Value = "'5.264"
$ws.Range("E4").Value = "'4.00%"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.07156"
$ws.Range("E5").Value = "'7.50%"
$ws.Range("G5").Value = "'22"
$ws.Range("D6").Value = "'7.516"
$ws.Range("E6").Value = "'2.08%"
$ws.Range("G6").Value = "'22"
$ws.Range("D7").Value = "'3.576"
$ws.Range("E7").Value = "'5.46%"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'1.409"
$ws.Range("E8").Value = "'2.16%"
$ws.Range("G8").Value = "'22"
$ws.Range("D9").Value = "'0.9104"
$ws.Range("E9").Value = "'-3.35%"
$ws.Range("G9").Value = "'22"
$ws.Range("D10").Value = "'0.1624"
$ws.Range("E10").Value = "'3.76%"
$ws.Range("G10").Value = "'22"
$ws.Range("D11").Value = "'0.07620"
$ws.Range("E11").Value = "'15.33%"
$ws.Range("G11").Value = "'22"
$ws.Range("D12").Value = "'0.07724"
$ws.Range("E12").Value = "'1.91%"
$ws.Range("G12").Value = "'22"
$ws.Range("D13").Value = "'0.02909"
$ws.Range("E13").Value = "'-1.09%"
$ws.Range("G13").Value = "'22"
$ws.Range("D14").Value = "'0.08999"
$ws.Range("E14").Value = "'0.02%"
$ws.Range("G14").Value = "'22"
$ws.Range("D15").Value = "'0.001588"
$ws.Range("E15").Value = "'-0.66%"
$ws.Range("G15").Value = "'22"
$ws.Range("D16").Value = "'0.0006527"
$ws.Range("E16").Value = "'1.41%"
$ws.Range("G16").Value = "'22"
$ws.Range("D17").Value = "'0.006257"
$ws.Range("E17").Value = "'-1.07%"
$ws.Range("G17").Value = "'22"
$ws.Range("D18").Value = "'3.475"
$ws.Range("E18").Value = "'0.58%"
$ws.Range("G18").Value = "'22"
$ws.Range("E19").Value = "'-1.24%"
$ws.Range("G19").Value = "'22"
$ws.Range("D20").Value = "'0.3269"
$ws.Range("E20").Value = "'1.69%"
$ws.Range("G20").Value = "'22"
$ws.Range("E21").Value = "'4.99%"
$ws.Range("G21").Value = "'22"
$ws.Range("D22").Value = "'4.017"
$ws.Range("E22").Value = "'-1.66%"
$ws.Range("G22").Value = "'22"
$ws.Range("E23").Value = "'2.51%"
$ws.Range("G23").Value = "'22"
$ws.Range("D24").Value = "'0.04523"
$ws.Range("E24").Value = "'1.04%"
$ws.Range("G24").Value = "'22"
$ws.Range("D25").Value = "'0.001205"
$ws.Range("G25").Value = "'22"
$ws.Range("D26").Value = "'0.004250"
$ws.Range("E26").Value = "'-5.32%"
$ws.Range("G26").Value = "'22"
$ws.Range("D27").Value = "'0.0001163"
$ws.Range("E27").Value = "'-7.02%"
$ws.Range("G27").Value = "'22"
$ws.Range("D28").Value = "'0.0001680"
$ws.Range("E28").Value = "'3.91%"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("G38").Value = "'22"
$ws.Range("G39").Value = "'22"
$ws.Range("D40").Value = "'0.04431"
$ws.Range("E40").Value = "'5.26%"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.007003"
$ws.Range("E41").Value = "'4.09%"
$ws.Range("G41").Value = "'22"
$ws.Range("D42").Value = "'0.1274"
$ws.Range("E42").Value = "'1.58%"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.002197"
$ws.Range("E43").Value = "'8.68%"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.01331"
$ws.Range("E44").Value = "'8.37%"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.00005809"
$ws.Range("E45").Value = "'2.17%"
$ws.Range("G45").Value = "'22"
$ws.Range("G46").Value = "'22"
$ws.Range("D47").Value = "'0.01292"
$ws.Range("E47").Value = "'-1.05%"
$ws.Range("G47").Value = "'22"
$ws.Range("G48").Value = "'22"
$ws.Range("G49").Value = "'22"
$ws.Range("G50").Value = "'22"
$ws.Range("G51").Value = "'22"
